$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (right after the header row), shifting all
# existing data rows (and their formatting) down by one, just like the
# daily price-scraper job does when it prepends today's entry.
$ws.Rows.Item(2).Insert()

# The inserted row picks up the bold header formatting by default; reset it
# back to the plain (unstyled) look used by every other data row.
$ws.Range("A2:D2").ClearFormats()

# Populate the newly inserted row with the latest price data entry. The
# leading apostrophe forces the date column to be stored as plain text
# (matching every other row, which holds its date as literal text rather
# than a real date value) instead of being auto-converted to a date serial.
$ws.Cells.Item(2, 1).Value = "'2025-12-24"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
